$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: clear the old "170/243" / Neural Network test entry ---
# A2, B2, E2 had no special style -> clear completely
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("E2").ClearContents()
# C2, D2 keep their percent-format style but lose their values
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

# --- New data rows 10-12 (entered first so shared-string order matches) ---
$ws.Range("A10").Value = "Arial"
$ws.Range("B10").Value = "144/167"
$ws.Range("C10").NumberFormat = "0.00%"
$ws.Range("C10").Value = 0.862

$ws.Range("A11").Value = "Arial"
$ws.Range("B11").Value = "159/167"
$ws.Range("C11").NumberFormat = "0.00%"
$ws.Range("C11").Value = 0.952

$ws.Range("A12").Value = "Arial"
$ws.Range("B12").Value = "208/243"
$ws.Range("C12").NumberFormat = "0.00%"
$ws.Range("C12").Value = 0.856

# --- New section header row 9: "AFTER OPTIMIZATION (Size invariant)" ---
$ws.Range("A9").Value = "AFTER OPTIMIZATION (Size invariant)"
$hdr = $ws.Range("A9:E9")
$hdr.Style = "Good"
$hdr.HorizontalAlignment = -4108
$hdr.Merge()

# --- Footnote on row 12 (entered last so shared-string order matches) ---
$ws.Range("F12").Value = '*** 212/243 after factoring in unavoidable OCR mistakes, such as the difference between "I" and "l" (lower L and Upper i)'

# --- Row 17: extend average formula to include new rows, apply percent style ---
$ws.Range("C17").NumberFormat = "0.00%"
$ws.Range("C17").Formula = "=AVERAGE(C2:C12)"

# --- Recalculate so cached formula values are correct ---
$excel.Calculate()

# --- Update the active selection as recorded in the saved file ---
$ws.Range("G16").Select()
